$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

# Update the conversion text on Hoja1!A1
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.13 = 20441.48 pesos`n✅ 20441.48 pesos = 5.11 = 932.4 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Update tasas sheet values
$ws2.Range("O10").Value = 3982
$ws2.Range("N12").Value = 3999.99
$ws2.Range("O12").Value = 182.453
